$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 122000
$ws.Range("E8").Value = 117800
$ws.Range("F8").Value = 151000
$ws.Range("G8").Value = 154000
$ws.Range("H8").Value = 138500
$ws.Range("I8").Value = 214400
$ws.Range("J8").Value = 221400
$ws.Range("D9").Value = 114500
$ws.Range("E9").Value = 122300
$ws.Range("F9").Value = 132300
$ws.Range("G9").Value = 138500
$ws.Range("H9").Value = 130400
$ws.Range("I9").Value = 157900
$ws.Range("J9").Value = 152000
$ws.Range("D10").Value = 7500
$ws.Range("E10").Value = -4500
$ws.Range("F10").Value = 18600
$ws.Range("G10").Value = 15500
$ws.Range("H10").Value = 8000
$ws.Range("I10").Value = 56500
$ws.Range("J10").Value = 69300
$ws.Range("D14").Value = 5400
$ws.Range("E14").Value = 34200
$ws.Range("F14").Value = 62600
$ws.Range("D17").Value = 133600
$ws.Range("E17").Value = 164700
$ws.Range("F17").Value = 201300
$ws.Range("G17").Value = 156100
$ws.Range("H17").Value = 139000
$ws.Range("I17").Value = 165600
$ws.Range("J17").Value = 161900
$ws.Range("D18").Value = -11600
$ws.Range("E18").Value = -46900
$ws.Range("F18").Value = -50300
$ws.Range("G18").Value = -2100
$ws.Range("I18").Value = 48800
$ws.Range("J18").Value = 59400
$ws.Range("D21").Value = -9300
$ws.Range("E21").Value = -40200
$ws.Range("F21").Value = -40100
$ws.Range("G21").Value = 8400
$ws.Range("H21").Value = 9900
$ws.Range("I21").Value = 58700
$ws.Range("J21").Value = 66700
$ws.Range("D23").Value = -11600
$ws.Range("E23").Value = -46900
$ws.Range("F23").Value = -50300
$ws.Range("G23").Value = -2100
$ws.Range("I23").Value = 48800
$ws.Range("J23").Value = 59400
$ws.Range("F24").Value = 3500
$ws.Range("I24").Value = 12600
$ws.Range("J24").Value = 15700
$ws.Range("D26").Value = -13100
$ws.Range("E26").Value = -47800
$ws.Range("F26").Value = -53800
$ws.Range("G26").Value = -4500
$ws.Range("I26").Value = 36200
$ws.Range("J26").Value = 43700
$ws.Range("D27").Value = -13100
$ws.Range("E27").Value = -47800
$ws.Range("F27").Value = -53800
$ws.Range("G27").Value = -4500
$ws.Range("I27").Value = 36200
$ws.Range("J27").Value = 43700
$ws.Range("D33").Value = -13100
$ws.Range("E33").Value = -47800
$ws.Range("F33").Value = -53800
$ws.Range("G33").Value = -4500
$ws.Range("I33").Value = 36200
$ws.Range("J33").Value = 43700
$ws.Range("D35").Value = -13100
$ws.Range("E35").Value = -47800
$ws.Range("F35").Value = -53800
$ws.Range("G35").Value = -4500
$ws.Range("I35").Value = 36200
$ws.Range("J35").Value = 43700
$ws.Range("G41").Value = 9100
$ws.Range("H41").Value = 4300
$ws.Range("I41").Value = 13300
$ws.Range("J41").Value = 6300
$ws.Range("D43").Value = 79300
$ws.Range("E43").Value = 84400
$ws.Range("F43").Value = 78500
$ws.Range("G43").Value = 81900
$ws.Range("H43").Value = 74600
$ws.Range("I43").Value = 67700
$ws.Range("J43").Value = 72900
$ws.Range("D44").Value = 28400
$ws.Range("E44").Value = 31600
$ws.Range("F44").Value = 45500
$ws.Range("G44").Value = 49100
$ws.Range("H44").Value = 45600
$ws.Range("I44").Value = 43100
$ws.Range("J44").Value = 43300
$ws.Range("F45").Value = 6200
$ws.Range("G45").Value = 4300
$ws.Range("H45").Value = 6100
$ws.Range("D46").Value = 108100
$ws.Range("E46").Value = 116000
$ws.Range("F46").Value = 130300
$ws.Range("G46").Value = 144400
$ws.Range("H46").Value = 130600
$ws.Range("I46").Value = 124700
$ws.Range("J46").Value = 123700
$ws.Range("D48").Value = 13700
$ws.Range("E48").Value = 20400
$ws.Range("F48").Value = 60800
$ws.Range("G48").Value = 106600
$ws.Range("H48").Value = 119100
$ws.Range("I48").Value = 118100
$ws.Range("J48").Value = 108800
$ws.Range("G49").Value = 4900
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5100
$ws.Range("J49").Value = 5200
$ws.Range("G52").Value = 1400
$ws.Range("H52").Value = 1500
$ws.Range("D54").Value = 122500
$ws.Range("E54").Value = 138200
$ws.Range("F54").Value = 194900
$ws.Range("G54").Value = 257300
$ws.Range("H54").Value = 256200
$ws.Range("I54").Value = 248100
$ws.Range("J54").Value = 254100
$ws.Range("D57").Value = 9100
$ws.Range("E57").Value = 12500
$ws.Range("F57").Value = 16400
$ws.Range("G57").Value = 18100
$ws.Range("H57").Value = 22600
$ws.Range("I57").Value = 17100
$ws.Range("J57").Value = 37500
$ws.Range("F58").Value = 5900
$ws.Range("G58").Value = 12600
$ws.Range("H58").Value = 14800
$ws.Range("J58").Value = 18600
$ws.Range("D59").Value = 10300
$ws.Range("E59").Value = 11100
$ws.Range("F59").Value = 11000
$ws.Range("G59").Value = 11500
$ws.Range("H59").Value = 7900
$ws.Range("I59").Value = 7700
$ws.Range("J59").Value = 15300
$ws.Range("D60").Value = 19400
$ws.Range("E60").Value = 23600
$ws.Range("F60").Value = 33400
$ws.Range("G60").Value = 42200
$ws.Range("H60").Value = 45300
$ws.Range("I60").Value = 26300
$ws.Range("J60").Value = 67800
$ws.Range("I61").Value = 7400
$ws.Range("J61").Value = 8900
$ws.Range("D66").Value = 19400
$ws.Range("E66").Value = 23600
$ws.Range("F66").Value = 33600
$ws.Range("G66").Value = 42400
$ws.Range("H66").Value = 45500
$ws.Range("I66").Value = 33900
$ws.Range("J66").Value = 76900
$ws.Range("D72").Value = 103100
$ws.Range("E72").Value = 114600
$ws.Range("F72").Value = 161300
$ws.Range("G72").Value = 214800
$ws.Range("H72").Value = 210700
$ws.Range("I72").Value = 214200
$ws.Range("J72").Value = 177100
$ws.Range("D76").Value = 103100
$ws.Range("E76").Value = 114600
$ws.Range("F76").Value = 161300
$ws.Range("G76").Value = 214800
$ws.Range("H76").Value = 210700
$ws.Range("I76").Value = 214200
$ws.Range("J76").Value = 177200
$ws.Range("D81").Value = -13100
$ws.Range("E81").Value = -47800
$ws.Range("F81").Value = -53800
$ws.Range("G81").Value = -4500
$ws.Range("I81").Value = 36200
$ws.Range("J81").Value = 43700
$ws.Range("E83").Value = 6700
$ws.Range("F83").Value = 10200
$ws.Range("G83").Value = 10400
$ws.Range("H83").Value = 10400
$ws.Range("I83").Value = 9900
$ws.Range("J83").Value = 7300
$ws.Range("F89").Value = 24000
$ws.Range("G89").Value = 3400
$ws.Range("H89").Value = 7700
$ws.Range("I89").Value = 28300
$ws.Range("J89").Value = 18400
$ws.Range("F91").Value = -24400
$ws.Range("H91").Value = -14400
$ws.Range("I91").Value = -2900
$ws.Range("J91").Value = -65500
$ws.Range("E94").Value = 6600
$ws.Range("F94").Value = -26200
$ws.Range("H94").Value = -19100
$ws.Range("J94").Value = -64300
$ws.Range("H96").Value = -1800
$ws.Range("E100").Value = -5000
$ws.Range("F100").Value = -7100
$ws.Range("H100").Value = 2400
$ws.Range("I100").Value = -18600
$ws.Range("J100").Value = 13100
$ws.Range("F102").Value = -9000
$ws.Range("G102").Value = 4800
$ws.Range("H102").Value = -9000
$ws.Range("I102").Value = 7000
$ws.Range("J102").Value = -32800
